$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("funding")

$ws.Range("B2").Value = "NA"
$ws.Range("C2").Value = "NA"
$ws.Range("D2").Value = "NA"

$ws.Activate()
$ws.Range("F11").Select()
